$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 data (date 2020-05-03 / serial 43954 already present in A53)
$ws.Range("B53").Value = 272
$ws.Range("C53").Value = 6465
$ws.Range("D53").Value = 1562
$ws.Range("E53").Value = 429
$ws.Range("F53").Value = 14
$ws.Range("G53").Value = 2041

# Row 54 data (date 2020-05-04 / serial 43955 already present in A54)
$ws.Range("B54").Value = 348
$ws.Range("C54").Value = 6813
$ws.Range("D54").Value = 1632
$ws.Range("E54").Value = 436
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 2139

# Update the active selection to reflect where the user ended up (G54)
$ws.Range("G54").Select()
